# Modifications pour utiliser XGBClassifier et ajuster les predictions
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Valeurs réelles" -> rename S+1/S+2/S+3 headers to *_class and
# replace the forecast probability columns (C,D,E) with classifier outputs.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Valeurs réelles")

$ws1.Range("C1").Value = "PRIX EXP POMME GALA FRANCE 136/200G CAT.I CAISSE_S+1_class"
$ws1.Range("D1").Value = "PRIX EXP POMME GALA FRANCE 136/200G CAT.I CAISSE_S+2_class"
$ws1.Range("E1").Value = "PRIX EXP POMME GALA FRANCE 136/200G CAT.I CAISSE_S+3_class"

$sheet1Values = @{
    "C2" = 3;  "D2" = 2;  "E2" = 1;
    "C3" = 2;  "D3" = 2;  "E3" = 2;
    "C4" = 2;  "D4" = 2;  "E4" = 2;
    "C5" = 2;  "D5" = 2;  "E5" = 2;
    "C6" = 2;  "D6" = 2;  "E6" = 2;
    "C7" = 2;  "D7" = 2;  "E7" = 0;
    "C8" = 2;  "D8" = 0;  "E8" = 4;
    "C9" = 0;  "D9" = 4;  "E9" = 4;
    "C10" = 4; "D10" = 4; "E10" = 1;
    "C11" = 4; "D11" = 1; "E11" = 1;
    "C12" = 1; "D12" = 1; "E12" = 4;
    "C13" = 1; "D13" = 4; "E13" = 2;
    "C14" = 4; "D14" = 2; "E14" = 4;
    "C15" = 2; "D15" = 4; "E15" = 2;
    "C16" = 4; "D16" = 2; "E16" = 0;
    "C17" = 2; "D17" = 0; "E17" = 1;
    "C18" = 0; "D18" = 1; "E18" = 3;
    "C19" = 1; "D19" = 3; "E19" = 0;
    "C20" = 3; "D20" = 0; "E20" = 1;
    "C21" = 0; "D21" = 1; "E21" = 3;
    "C22" = 1; "D22" = 3; "E22" = 1;
    "C23" = 3; "D23" = 1; "E23" = 3;
    "C24" = 1; "D24" = 3; "E24" = 4;
    "C25" = 3; "D25" = 4; "E25" = 4;
    "C26" = 4; "D26" = 4; "E26" = 2;
    "C27" = 4; "D27" = 2; "E27" = 2;
    "C28" = 2; "D28" = 2; "E28" = 2;
}

foreach ($addr in $sheet1Values.Keys) {
    $ws1.Range($addr).Value = $sheet1Values[$addr]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Prédictions" -> replace probability predictions (B,C,D) with
# classifier-based integer predictions.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Prédictions")

$sheet2Values = @{
    "B2" = 1;  "C2" = 0;  "D2" = 0;
    "B3" = 1;  "C3" = -1; "D3" = 0;
    "B4" = 1;  "C4" = -1; "D4" = 0;
    "B5" = 1;  "C5" = -1; "D5" = 0;
    "B6" = 1;  "C6" = -1; "D6" = 0;
    "B7" = 1;  "C7" = -1; "D7" = 0;
    "B8" = 1;  "C8" = -1; "D8" = 0;
    "B9" = 1;  "C9" = -1; "D9" = 0;
    "B10" = 2; "C10" = 0;  "D10" = 0;
    "B11" = 1; "C11" = -1; "D11" = 0;
    "B12" = 1; "C12" = -1; "D12" = 0;
    "B13" = 0; "C13" = 0;  "D13" = 0;
    "B14" = 1; "C14" = 0;  "D14" = 0;
    "B15" = 0; "C15" = 0;  "D15" = 0;
    "B16" = 0; "C16" = 0;  "D16" = 0;
    "B17" = 1; "C17" = 0;  "D17" = 0;
    "B18" = 0; "C18" = 0;  "D18" = 0;
    "B19" = 0; "C19" = -1; "D19" = 0;
    "B20" = 1; "C20" = 0;  "D20" = 0;
    "B21" = 0; "C21" = 0;  "D21" = 0;
    "B22" = 0; "C22" = 0;  "D22" = 0;
    "B23" = 0; "C23" = 0;  "D23" = 0;
    "B24" = 1; "C24" = 0;  "D24" = 0;
    "B25" = 0; "C25" = 1;  "D25" = 0;
    "B26" = -1; "C26" = 0; "D26" = 0;
    "B27" = -1; "C27" = 1; "D27" = 1;
    "B28" = 0; "C28" = 1;  "D28" = 0;
}

foreach ($addr in $sheet2Values.Keys) {
    $ws2.Range($addr).Value = $sheet2Values[$addr]
}
